$d = $word.ActiveDocument

# 1. Remove the "_GoBack" bookmark that sits in the "2-3" resources table cell.
#    (It will be re-created later near the ethics paragraph, see step 4.)
foreach ($bm in @($d.Bookmarks)) {
    if ($bm.Name -eq "_GoBack") {
        $bm.Delete()
    }
}

# 2. Split "... incremental approach with x increments:" so the "x" becomes "7"
#    as its own run.
$d.Content.Find.Execute("incremental approach with x increments:", $true, $false, $false, $false, $false, $true, 1, $false, "incremental approach with 7 increments:", 2)

# 3. "I intend to stick within..." -> "The project aims to stick within..."
$d.Content.Find.Execute("I intend to stick", $true, $false, $false, $false, $false, $true, 1, $false, "The project aims to stick", 2)

# 4. Re-add the _GoBack bookmark right after the ethics-application sentence,
#    and remove the trailing empty paragraph before the section break.
$rng = $d.Content
$rng.Find.Execute("stick within the limitations of the approved ethics application.")
$endOfSentence = $rng.End
$d.Bookmarks.Add("_GoBack", $d.Range($endOfSentence, $endOfSentence))

$paras = $d.Paragraphs
$lastPara = $paras.Item($paras.Count)
if ($lastPara.Range.Text -eq "`r") {
    $lastPara.Range.Delete()
}
